# Auto-generated edit script: apply scheduled market-data refresh to Sheets
# Updates cached currentAveragePrice / LevePrice / LeveProfit columns (H:N)
# per leve row across all 8 job sheets. Values are plain cached numbers
# (no formulas in this workbook), so we write literals directly.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 809.2
$ws.Range("I39").Value = 809.2
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 2427.6
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = -2131.6
$ws.Range("N39").ClearContents()

$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()

$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()

$ws.Range("H86").Value = 4247.923
$ws.Range("I86").Value = 1887.25
$ws.Range("K86").Value = 1887.25
$ws.Range("M86").Value = -764.25

$ws.Range("H88").Value = 6078.091
$ws.Range("I88").Value = 791.6
$ws.Range("J88").Value = 10483.5
$ws.Range("K88").Value = 791.6
$ws.Range("L88").Value = 10483.5
$ws.Range("M88").Value = -385.6
$ws.Range("N88").Value = -11295.5

$ws.Range("H89").Value = 4247.923
$ws.Range("I89").Value = 1887.25
$ws.Range("K89").Value = 9436.25
$ws.Range("M89").Value = -3820.25

$ws.Range("H91").Value = 6078.091
$ws.Range("I91").Value = 791.6
$ws.Range("J91").Value = 10483.5
$ws.Range("K91").Value = 791.6
$ws.Range("L91").Value = 10483.5
$ws.Range("M91").Value = 612.4
$ws.Range("N91").Value = -13291.5

$ws.Range("H99").Value = 1412.8334
$ws.Range("J99").Value = 2393.6667
$ws.Range("L99").Value = 7181.000100000001
$ws.Range("N99").Value = -10177.0001

$ws.Range("H107").Value = 505.9
$ws.Range("I107").Value = 504.85715
$ws.Range("J107").Value = 508.33334
$ws.Range("K107").Value = 504.85715
$ws.Range("L107").Value = 508.33334
$ws.Range("M107").Value = 1415.14285
$ws.Range("N107").Value = -4348.33334

$ws.Range("H112").Value = 3199.5557
$ws.Range("I112").Value = 2149.5
$ws.Range("J112").Value = 3499.5715
$ws.Range("K112").Value = 6448.5
$ws.Range("L112").Value = 10498.7145
$ws.Range("M112").Value = -5340.5
$ws.Range("N112").Value = -12714.7145

$ws.Range("H129").Value = 2427.7144
$ws.Range("I129").Value = 2000
$ws.Range("J129").Value = 2460.6155
$ws.Range("K129").Value = 6000
$ws.Range("L129").Value = 7381.8465
$ws.Range("M129").Value = -1000
$ws.Range("N129").Value = -17381.8465

$ws.Range("H132").Value = 934.5
$ws.Range("I132").Value = 872
$ws.Range("K132").Value = 2616
$ws.Range("M132").Value = -86

$ws.Range("H138").Value = 4078.2856
$ws.Range("I138").Value = 2580.1904
$ws.Range("J138").Value = 5201.857
$ws.Range("K138").Value = 7740.5712
$ws.Range("L138").Value = 15605.571
$ws.Range("M138").Value = -2600.5712
$ws.Range("N138").Value = -25885.571

$ws.Range("H141").Value = 2012.5714
$ws.Range("I141").Value = 2012.5714
$ws.Range("K141").Value = 6037.7142
$ws.Range("M141").Value = -857.7142000000003


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4684.8965
$ws.Range("I32").Value = 4495.143
$ws.Range("J32").Value = 9998
$ws.Range("K32").Value = 4495.143
$ws.Range("L32").Value = 9998
$ws.Range("M32").Value = -4208.143
$ws.Range("N32").Value = -10572

$ws.Range("H138").Value = 1650000
$ws.Range("J138").Value = 1650000
$ws.Range("L138").Value = 1650000
$ws.Range("N138").Value = -1660280


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 8148.8887
$ws.Range("I94").Value = 8148.8887
$ws.Range("K94").Value = 8148.8887
$ws.Range("M94").Value = -7697.8887

$ws.Range("H105").Value = 4000
$ws.Range("I105").Value = 4000
$ws.Range("K105").Value = 4000
$ws.Range("M105").Value = -2253

$ws.Range("H134").Value = 1372.6
$ws.Range("I134").Value = 1372.6
$ws.Range("K134").Value = 4117.799999999999
$ws.Range("M134").Value = -1582.799999999999


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 25443.7
$ws.Range("I22").Value = 554.125
$ws.Range("K22").Value = 554.125
$ws.Range("M22").Value = -204.125


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 12500115
$ws.Range("I4").Value = 12500115
$ws.Range("K4").Value = 37500345
$ws.Range("M4").Value = -37500233

$ws.Range("H5").Value = 5399
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 5399
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 16197
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -16421

$ws.Range("H11").Value = 25017862
$ws.Range("I11").Value = 25017862
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 75053586
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -75053446
$ws.Range("N11").ClearContents()

$ws.Range("H29").Value = 1630
$ws.Range("J29").Value = 2150
$ws.Range("L29").Value = 6450
$ws.Range("N29").Value = -7004

$ws.Range("H39").Value = 8412
$ws.Range("J39").Value = 8412
$ws.Range("L39").Value = 25236
$ws.Range("N39").Value = -25824

$ws.Range("H55").Value = 4601
$ws.Range("J55").Value = 6668.3335
$ws.Range("L55").Value = 20005.0005
$ws.Range("N55").Value = -20359.0005

$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("M81").ClearContents()

$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("M84").ClearContents()

$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("M94").ClearContents()
$ws.Range("N94").ClearContents()

$ws.Range("H113").Value = 1085.6364
$ws.Range("J113").Value = 995.3333
$ws.Range("L113").Value = 2985.9999
$ws.Range("N113").Value = -7325.9999

$ws.Range("H115").Value = 0
$ws.Range("I115").Value = 0
$ws.Range("K115").Value = 0
$ws.Range("M115").ClearContents()

$ws.Range("H118").Value = 4999
$ws.Range("I118").Value = 4999
$ws.Range("K118").Value = 14997
$ws.Range("M118").Value = -13754

$ws.Range("H132").Value = 4295.7
$ws.Range("I132").Value = 4279.7144
$ws.Range("J132").Value = 4333
$ws.Range("K132").Value = 38517.4296
$ws.Range("L132").Value = 38997
$ws.Range("M132").Value = -35987.4296
$ws.Range("N132").Value = -44057

$ws.Range("H135").Value = 5399
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 5399
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 48591
$ws.Range("M135").ClearContents()
$ws.Range("N135").Value = -53661

$ws.Range("H137").Value = 2332.3333
$ws.Range("J137").Value = 3622.75
$ws.Range("L137").Value = 10868.25
$ws.Range("N137").Value = -21068.25

$ws.Range("H139").Value = 1392.6
$ws.Range("I139").Value = 1392.6
$ws.Range("K139").Value = 4177.799999999999
$ws.Range("M139").Value = 962.2000000000007


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 7.5
$ws.Range("I4").Value = 5
$ws.Range("J4").Value = 10
$ws.Range("K4").Value = 5
$ws.Range("L4").Value = 10
$ws.Range("M4").Value = 107
$ws.Range("N4").Value = -234

$ws.Range("H5").Value = 4.3333335
$ws.Range("I5").Value = 4.3333335
$ws.Range("K5").Value = 4.3333335
$ws.Range("M5").Value = 107.6666665

$ws.Range("H92").Value = 50251
$ws.Range("J92").Value = 50251
$ws.Range("L92").Value = 50251
$ws.Range("N92").Value = -53995


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 751.5833
$ws.Range("J16").Value = 1789.6666
$ws.Range("L16").Value = 1789.6666
$ws.Range("N16").Value = -2129.6666

$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").ClearContents()

$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("N27").ClearContents()

$ws.Range("H82").Value = 1318.909
$ws.Range("I82").Value = 1378.3334
$ws.Range("K82").Value = 1378.3334
$ws.Range("M82").Value = -1017.3334

$ws.Range("H85").Value = 1318.909
$ws.Range("I85").Value = 1378.3334
$ws.Range("K85").Value = 1378.3334
$ws.Range("M85").Value = -130.3334

$ws.Range("H122").Value = 7267.607
$ws.Range("I122").Value = 7623.4707
$ws.Range("K122").Value = 22870.4121
$ws.Range("M122").Value = -20420.4121

$ws.Range("H136").Value = 2619.8462
$ws.Range("I136").Value = 2192.7646
$ws.Range("K136").Value = 6578.293799999999
$ws.Range("M136").Value = -4028.293799999999


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 12608.23
$ws.Range("I62").Value = 8793.6
$ws.Range("J62").Value = 14992.375
$ws.Range("K62").Value = 8793.6
$ws.Range("L62").Value = 14992.375
$ws.Range("M62").Value = -8169.6
$ws.Range("N62").Value = -16240.375

$ws.Range("H65").Value = 12608.23
$ws.Range("I65").Value = 8793.6
$ws.Range("J65").Value = 14992.375
$ws.Range("K65").Value = 43968
$ws.Range("L65").Value = 74961.875
$ws.Range("M65").Value = -40848
$ws.Range("N65").Value = -81201.875

$ws.Range("H107").Value = 2000
$ws.Range("I107").Value = 2000
$ws.Range("K107").Value = 6000
$ws.Range("M107").Value = -4080

$ws.Range("H132").Value = 2187.6667
$ws.Range("I132").Value = 555.7143
$ws.Range("J132").Value = 7899.5
$ws.Range("K132").Value = 1667.1429
$ws.Range("L132").Value = 23698.5
$ws.Range("M132").Value = 862.8571000000002
$ws.Range("N132").Value = -28758.5

$ws.Range("H136").Value = 2145.6155
$ws.Range("I136").Value = 2157.75
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 6473.25
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -3923.25
$ws.Range("N136").Value = -11100

